$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item(1)
$wsZhCn     = $wb.Worksheets.Item(2)
$wsDeDe     = $wb.Worksheets.Item(3)

# Update the status text from "Ready for handoff" to "In Translation"
# everywhere it is used (Overview!E2, Overview!F2, zh-cn!C2, de-de!C2).
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value     = "In Translation"
$wsDeDe.Range("C2").Value     = "In Translation"

# Shrink the columns that held the (now shorter) status text to their
# new auto-fitted width.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth     = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth     = 12.5
